$wb = $excel.ActiveWorkbook

# --- Sheet "Reguły": reorder the object lists inside a few rule descriptions ---
$rules = $wb.Worksheets.Item("Reguły")

$rules.Range("B2").Value = "(attempts >=  3.0) & (pregnancy <=  0.0) => (class <= 1) ['a1', 'a3', 'a7']"
$rules.Range("B6").Value = "(age >=  42.0) => (class <= 1) ['a3', 'a14']"
$rules.Range("B7").Value = "(age <=  31.0) & (attempts <=  1.0) & (endometrium <=  1.0) => (class >= 2) ['a11', 'a24', 'a12', 'a9']"
$rules.Range("B8").Value = "(frozen_embryos >=  8.0) & (sperm <=  1.0) => (class >= 2) ['a6', 'a16']"

# --- Sheet "Walidacja krzyżowa": swap the "accuracy" and "correct" rows (label + value) ---
$cv = $wb.Worksheets.Item("Walidacja krzyżowa")

$cv.Range("A1").Value = "correct"
$cv.Range("B1").Value = 0.8181818181818182
$cv.Range("A4").Value = "accuracy"
$cv.Range("B4").Value = 0.36
